$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), formatted the same as the
# existing header cells (bold, centered, thin-bordered -> style index 1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells for row 2 (plain, unstyled numbers).
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9
